# fixed path in excel files
# - source/target sql path strings now use forward slashes instead of backslashes
# - selection on the sheet moved from B8 to B10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the source/target SQL path values (row 15 = sourcequerysqlpath,
# row 28 = targetquerysqlpath) to use forward slashes.
$ws.Range("B15").Value = "test/sql/sourcesql"
$ws.Range("B28").Value = "test/sql/targetsql"

# Move the active selection from B8 to B10.
$ws.Range("B10").Select()
